$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Package tracking numbers in column C (and mirrored into column D for a
# subset of rows) are being refreshed with newer values. These look like
# numbers but must stay text (shared-string) cells, so the number format is
# forced to Text ("@") before assignment and reset back to Normal afterwards
# so the resulting cells carry no explicit style, matching the original
# formatting of these cells.

$textCells = $ws.Range("C2:C22,D5:D7,D13:D17")
$textCells.NumberFormat = "@"

$ws.Range("C2").Value = "320018407199"
$ws.Range("C3").Value = "320018407203"
$ws.Range("C4").Value = "320018407236"

$ws.Range("C5").Value = "320018407269"
$ws.Range("D5").Value = "320018407269"

$ws.Range("C6").Value = "320018407306"
$ws.Range("D6").Value = "320018407306"

$ws.Range("C7").Value = "320018407328"
$ws.Range("D7").Value = "320018407328"

$ws.Range("C8").Value = "320018407361"
$ws.Range("C9").Value = "320018407394"
$ws.Range("C10").Value = "320018407420"
$ws.Range("C11").Value = "320018407442"
$ws.Range("C12").Value = "320018407486"

$ws.Range("C13").Value = "320018407501"
$ws.Range("D13").Value = "320018407501"

$ws.Range("C14").Value = "320018407740"
$ws.Range("D14").Value = "320018407740"

$ws.Range("C15").Value = "320018407773"
$ws.Range("D15").Value = "320018407773"

$ws.Range("C16").Value = "320018407810"
$ws.Range("D16").Value = "320018407810"

$ws.Range("C17").Value = "320018407832"
$ws.Range("D17").Value = "320018407832"

$ws.Range("C18").Value = "320018407876"
$ws.Range("C19").Value = "320018407898"
$ws.Range("C20").Value = "320018407924"
$ws.Range("C21").Value = "320018407946"
$ws.Range("C22").Value = "320018407979"

$textCells.Style = "Normal"
